$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 38 (shifts old rows 38:43 down to 39:44)
$ws.Rows("38:38").Insert()

# Populate the newly inserted row 38 with the latest weekly entry
$ws.Range("A38").Value = 7
$ws.Range("B38").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C38").Value = "Ñuble"
$ws.Range("D38").Value = 44798
$ws.Range("E38").Value = 16
$ws.Range("F38").Value = 100112001
$ws.Range("G38").Value = "Berenjena"
$ws.Range("H38").Value = "Sin especificar"
$ws.Range("I38").Value = "Primera"
$ws.Range("J38").Value = 40
$ws.Range("K38").Value = 12000
$ws.Range("L38").Value = 12000
$ws.Range("M38").Value = 12000
$ws.Range("N38").Value = "`$/caja 60 unidades"
$ws.Range("O38").Value = "Región de Arica y Parinacota"
$ws.Range("P38").Value = 200
$ws.Range("Q38").Value = 60
$ws.Range("R38").Value = "Hortaliza"
